# Weekly update: insert a new price record for Espárragos at
# "Vega Modelo de Temuco" (row 17) and push the existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17, shifting rows 17:47
# down to 18:48 (and extending the used range to A1:R48).
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with this week's record.
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44498
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 300000000
$ws.Range("G17").Value = "Espárragos"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 1300
$ws.Range("M17").Value = 1250
$ws.Range("N17").Value = "$/kilo"
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 1250
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"
